# 6.3.1.xlsx — add the 2021 data column (O) to the right of the existing
# 2010-2020 columns (D-N), mirroring the formatting of column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column N's formatting (borders/fonts/number formats/etc.) into column
# O for the data block (header row through the last data row) so the new
# column renders identically to the others.
$ws.Range("N3:N14").Copy()
$ws.Range("O3:O14").PasteSpecial(-4122)

# Fill in the 2021 values.
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 97
$ws.Range("O6").Value = 96.2
$ws.Range("O7").Value = 62.7
$ws.Range("O8").Value = 100
$ws.Range("O9").Value = 100
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = 100
$ws.Range("O12").Value = 57.9
$ws.Range("O13").Value = 100
$ws.Range("O14").Value = "-"

# Match the saved selection from the edited workbook.
$ws.Range("O17").Select()
